# Generate Report for Handback
# - Update the "Status" column for the 0a359e45-... file (row 3) from
#   "Ready for handoff" to "Handback transform failed" on the Overview,
#   zh-cn and de-de sheets.
# - Populate the "Error Detail" column (L) for row 3 on the zh-cn and
#   de-de sheets with the handback/handoff file name mismatch message.

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

# --- Overview sheet: both locale status cells (B3 = zh-cn, C3 = de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("L3").Value = "Handback file name: nwibzvfb.j0q is different with handoff file name: 0a359e45-d336-4588-a6a7-a296bf4df260.517da1f4e64dc06ecb55192ce70ddcb1ea87c429.zh-cn."

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("L3").Value = "Handback file name: nwibzvfb.j0q is different with handoff file name: 0a359e45-d336-4588-a6a7-a296bf4df260.517da1f4e64dc06ecb55192ce70ddcb1ea87c429.de-de."
